# Updated to work on latest Dynamics version
#
# RVL sheet: after playing the PriceListItems sheet, also play the Cleanup
# sheet (closes the Navigator) before the flow ends.
#
# PriceListItems sheet: the newer Dynamics UI requires clicking into the
# "Pricing_information" section before the Amount field can be edited.

$wb = $excel.ActiveWorkbook

$rvl = $wb.Worksheets.Item("RVL")
$rvl.Range("B45").Value = "Action"
$rvl.Range("C45").Value = "RVL"
$rvl.Range("D45").Value = "DoPlaySheet"
$rvl.Range("E45").Value = "sheetName"
$rvl.Range("F45").Value = "string"
$rvl.Range("G45").Value = "Cleanup"

$rvl.Range("B46").Value = "Action"
$rvl.Range("C46").Value = "Navigator"
$rvl.Range("D46").Value = "Close"

$pli = $wb.Worksheets.Item("PriceListItems")
$pli.Rows.Item(26).Insert()
$pli.Range("B26").Value = "Action"
$pli.Range("C26").Value = "Pricing_information"
$pli.Range("D26").Value = "DoClick"
